$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: "HUS, PPE Applied Year" should be 2021, not 2016
$ws.Range("B2").Value = 2021

# Match the selection state observed in the saved file (active cell B2)
$ws.Range("B2").Select()

$wb.Save()
